# Update cryptocurrency price/volume data in the worksheet
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "67.480.99"
$ws.Range("E2").Value = "  -1.04%  "

# Row 3
$ws.Range("D3").Value = "3.224.82"
$ws.Range("E3").Value = "  -1.54%  "

# Row 4
$ws.Range("E4").Value = "  -0.03%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "578.24"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.63%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "183.39"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.30%  "

# Row 7
$ws.Range("E7").Value = "  -0.04%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.604"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.34%  "

# Row 9
$ws.Range("D9").Value = "3.223.34"
$ws.Range("E9").Value = "  -1.59%  "

# Row 10
$ws.Range("E10").Value = "  -3.17%  "

# Row 12
$ws.Range("E12").Value = "  -1.35%  "

# Row 13
$ws.Range("D13").Value = "3.777.29"
$ws.Range("E13").Value = "  -1.84%  "

# Row 14
$ws.Range("E14").Value = "  +0.10%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "27.78"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -3.39%  "

# Row 16
$ws.Range("D16").Value = "67.527.85"
$ws.Range("E16").Value = "  -0.99%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.0000169"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -2.32%  "

# Row 18
$ws.Range("D18").Value = "3.241.58"
$ws.Range("E18").Value = "  -0.81%  "

# Row 19
$ws.Range("E19").Value = "  -2.29%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "13.42"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.82%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "395.49"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +3.10%  "

# Row 22
$ws.Range("E22").Value = "  -2.33%  "

# Row 23
$ws.Range("E23").Value = "  +0.00%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "71.16"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.50%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.514"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.61%  "

# Row 26
$ws.Range("E26").Value = "  -2.96%  "

# Row 27
$ws.Range("E27").Value = "  -0.25%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.57"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -3.61%  "

# Row 29
$ws.Range("E29").Value = "  +0.07%  "

# Row 30
$ws.Range("E30").Value = "  -2.56%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "5.57"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -4.33%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "22.60"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -1.75%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "6.95"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -4.48%  "

# Row 35
$ws.Range("E35").Value = "  -2.91%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "160.87"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -1.11%  "

# Row 37
$ws.Range("E37").Value = "  -5.07%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.88"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.05%  "

# Row 39
$ws.Range("B39").Value = "EnergySwap"
$ws.Range("C39").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "26.22"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -2.38%  "

# Row 40
$ws.Range("B40").Value = "Mantle"
$ws.Range("C40").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.802"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -4.41%  "

# Row 41
$ws.Range("E41").Value = "  -1.60%  "

# Row 42
$ws.Range("E42").Value = "  -4.54%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.47"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -5.57%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0684"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.85%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "40.50"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -2.54%  "

# Row 46
$ws.Range("D46").Value = "2.590.44"
$ws.Range("E46").Value = "  -2.57%  "

# Row 47
$ws.Range("B47").Value = "InjectiveProtocol"
$ws.Range("C47").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "24.50"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -4.41%  "

# Row 48
$ws.Range("B48").Value = "Bittensor"
$ws.Range("C48").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "333.60"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -4.49%  "

# Row 49
$ws.Range("E49").Value = "  -2.82%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "6.28"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.02%  "

# Row 51
$ws.Range("E51").Value = "  -1.99%  "
